$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.770123734158688
$ws.Range("D2").Value = 0.003338351862899103
$ws.Range("E2").Value = 1.252180992107384
$ws.Range("F2").Value = 0.6385628759164916
$ws.Range("G2").Value = 0.5072477128966142
$ws.Range("H2").Value = 0.5500569983951777
$ws.Range("L2").Value = 0.855296197323014

$ws.Range("B3").Value = 1.642592827950239
$ws.Range("D3").Value = 0.00351852692437582
$ws.Range("E3").Value = 1.135848692201591
$ws.Range("F3").Value = 0.6196949980858193
$ws.Range("G3").Value = 0.485844276608745
$ws.Range("H3").Value = 0.5493356329887149
$ws.Range("L3").Value = 0.7525873294032408

$ws.Range("B4").Value = 1.564482838566448
$ws.Range("D4").Value = 0.003648686322375028
$ws.Range("E4").Value = 1.064405800995274
$ws.Range("F4").Value = 0.6091672056908948
$ws.Range("G4").Value = 0.4737199112432364
$ws.Range("H4").Value = 0.5496884397640542
$ws.Range("L4").Value = 0.6893685676346024

$ws.Range("B5").Value = 1.532701153963103
$ws.Range("D5").Value = 0.00370655126145536
$ws.Range("E5").Value = 1.035290613362577
$ws.Range("F5").Value = 0.6051391191049476
$ws.Range("G5").Value = 0.469030445300362
$ws.Range("H5").Value = 0.550030338678539
$ws.Range("L5").Value = 0.6635674818862469

$ws.Range("B6").Value = 1.527426772819695
$ws.Range("D6").Value = 0.003716448517431914
$ws.Range("E6").Value = 1.030456020840859
$ws.Range("F6").Value = 0.6044859773283804
$ws.Range("G6").Value = 0.468266809342893
$ws.Range("H6").Value = 0.5500990197551374
$ws.Range("L6").Value = 0.6592808839926079

$ws.Range("B7").Value = 1.564054022304447
$ws.Range("D7").Value = 0.003649447291616426
$ws.Range("E7").Value = 1.064013147276967
$ws.Range("F7").Value = 0.609111825330757
$ws.Range("G7").Value = 0.4736556558899565
$ws.Range("H7").Value = 0.5496922511095335
$ws.Range("L7").Value = 0.6890207626307188

$ws.Range("B8").Value = 1.726110558944356
$ws.Range("D8").Value = 0.003396369526026177
$ws.Range("E8").Value = 1.212073725216953
$ws.Range("F8").Value = 0.6318355462696132
$ws.Range("G8").Value = 0.499653903795263
$ws.Range("H8").Value = 0.5496419315555556
$ws.Range("L8").Value = 0.819914361445143

$ws.Range("B9").Value = 2.045473510088129
$ws.Range("D9").Value = 0.003059048911914886
$ws.Range("E9").Value = 1.502238295003508
$ws.Range("F9").Value = 0.6849573692318245
$ws.Range("G9").Value = 0.5589164174766097
$ws.Range("H9").Value = 0.5559467457455582
$ws.Range("L9").Value = 1.07537901880255

$ws.Range("B10").Value = 2.281135022391311
$ws.Range("D10").Value = 0.002913735497218184
$ws.Range("E10").Value = 1.715249659232029
$ws.Range("F10").Value = 0.7294482764373384
$ws.Range("G10").Value = 0.6077953319822882
$ws.Range("H10").Value = 0.5646085460074914
$ws.Range("L10").Value = 1.262369446990931

$ws.Range("B11").Value = 2.388584198396586
$ws.Range("D11").Value = 0.002871133208127929
$ws.Range("E11").Value = 1.812107136054436
$ws.Range("F11").Value = 0.7509274965451738
$ws.Range("G11").Value = 0.6312541892335162
$ws.Range("H11").Value = 0.5694514539569866
$ws.Range("L11").Value = 1.347296310821321

$ws.Range("B12").Value = 2.429308768186559
$ws.Range("D12").Value = 0.002858485361358021
$ws.Range("E12").Value = 1.848777318871896
$ws.Range("F12").Value = 0.7592437081703736
$ws.Range("G12").Value = 0.6403185096981758
$ws.Range("H12").Value = 0.5714173611673914
$ws.Range("L12").Value = 1.379437015373298

$ws.Range("B13").Value = 2.420536396474688
$ws.Range("D13").Value = 0.002861052617010529
$ws.Range("E13").Value = 1.840880101981497
$ws.Range("F13").Value = 0.7574444797370319
$ws.Range("G13").Value = 0.6383582176520406
$ws.Range("H13").Value = 0.5709880616189196
$ws.Range("L13").Value = 1.372515790244677

$ws.Range("B14").Value = 2.391933910107184
$ws.Range("D14").Value = 0.002870022212309919
$ws.Range("E14").Value = 1.815124175595827
$ws.Range("F14").Value = 0.7516079947284879
$ws.Range("G14").Value = 0.6319962607827563
$ws.Range("H14").Value = 0.5696105320679692
$ws.Range("L14").Value = 1.349940929797754

$ws.Range("B15").Value = 2.374418755409238
$ws.Range("D15").Value = 0.002875973323257952
$ws.Range("E15").Value = 1.799346888241871
$ws.Range("F15").Value = 0.7480568701491279
$ws.Range("G15").Value = 0.6281230942325351
$ws.Range("H15").Value = 0.5687840120207852
$ws.Range("L15").Value = 1.336110678750401

$ws.Range("B16").Value = 2.274117955499548
$ws.Range("D16").Value = 0.002917001511198336
$ws.Range("E16").Value = 1.708918806992159
$ws.Range("F16").Value = 0.7280698686031712
$ws.Range("G16").Value = 0.6062872896081046
$ws.Range("H16").Value = 0.5643104017932217
$ws.Range("L16").Value = 1.256816563696077

$ws.Range("B17").Value = 2.212649947278408
$ws.Range("D17").Value = 0.002948264429790015
$ws.Range("E17").Value = 1.653432103453326
$ws.Range("F17").Value = 0.7161290650273173
$ws.Range("G17").Value = 0.5932086740178875
$ws.Range("H17").Value = 0.5617987326854177
$ws.Range("L17").Value = 1.20813754402883

$ws.Range("B18").Value = 2.177318248141489
$ws.Range("D18").Value = 0.002968451677828554
$ws.Range("E18").Value = 1.621513706948264
$ws.Range("F18").Value = 0.7093775015155899
$ws.Range("G18").Value = 0.5858010688918682
$ws.Range("H18").Value = 0.5604389119495465
$ws.Range("L18").Value = 1.180125833093427

$ws.Range("B19").Value = 2.16535948847104
$ws.Range("D19").Value = 0.002975662384717026
$ws.Range("E19").Value = 1.610706063105965
$ws.Range("F19").Value = 0.7071114154572911
$ws.Range("G19").Value = 0.5833125578908493
$ws.Range("H19").Value = 0.5599930017761778
$ws.Range("L19").Value = 1.170639338278306

$ws.Range("B20").Value = 2.219190929732008
$ws.Range("D20").Value = 0.002944707407992908
$ws.Range("E20").Value = 1.659339174154184
$ws.Range("F20").Value = 0.7173880995323856
$ws.Range("G20").Value = 0.5945889903514967
$ws.Range("H20").Value = 0.5620573091601955
$ws.Range("L20").Value = 1.213320834314629

$ws.Range("B21").Value = 2.400334174265254
$ws.Range("D21").Value = 0.00286729220118076
$ws.Range("E21").Value = 1.822689534565967
$ws.Range("F21").Value = 0.7533173254880978
$ws.Range("G21").Value = 0.6338599703992429
$ws.Range("H21").Value = 0.570011545663192
$ws.Range("L21").Value = 1.356572233287409

$ws.Range("B22").Value = 2.518931340754193
$ws.Range("D22").Value = 0.00283705893572872
$ws.Range("E22").Value = 1.929403266300113
$ws.Range("F22").Value = 0.7778648475484431
$ws.Range("G22").Value = 0.6605829552904368
$ws.Range("H22").Value = 0.5759805531567395
$ws.Range("L22").Value = 1.450083267988987

$ws.Range("B23").Value = 2.455614405876076
$ws.Range("D23").Value = 0.002851296507301981
$ws.Range("E23").Value = 1.872452750045284
$ws.Range("F23").Value = 0.7646644731389927
$ws.Range("G23").Value = 0.6462220022613394
$ws.Range("H23").Value = 0.5727235499865344
$ws.Range("L23").Value = 1.400184797893473

$ws.Range("B24").Value = 2.216233729562248
$ws.Range("D24").Value = 0.002946308650933815
$ws.Range("E24").Value = 1.656668644532886
$ws.Range("F24").Value = 0.7168185372460272
$ws.Range("G24").Value = 0.5939646025265688
$ws.Range("H24").Value = 0.5619401447281973
$ws.Range("L24").Value = 1.210977548206415

$ws.Range("B25").Value = 1.958902541013515
$ws.Range("D25").Value = 0.003132713055075342
$ws.Range("E25").Value = 1.423767370630344
$ws.Range("F25").Value = 0.6696449648022309
$ws.Range("G25").Value = 0.5419679397720358
$ws.Range("H25").Value = 0.5535426311665788
$ws.Range("L25").Value = 1.006395285445933
